# Rename a few transcript speaker labels in the "Speaker" column (D).
# "RBD" -> "T", "Students" -> "SS", "Student" -> "S"
# (trailing punctuation / whitespace in the original text is preserved)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2

    if ($val -eq $null) { continue }

    if ($val -eq "RBD") {
        $cell.Value = "T"
    }
    elseif ($val -eq "Students") {
        $cell.Value = "SS"
    }
    elseif ($val -eq "Students. ") {
        $cell.Value = "SS. "
    }
    elseif ($val -eq "Student") {
        $cell.Value = "S"
    }
}
